$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 6 (pushes existing rows 6-16 down to 8-18),
# copying formatting down from the row above just like Excel's native
# "Insert Sheet Rows" does.
$ws.Rows("6:7").Insert()

# --- Row 6: new timesheet entry (Doyle, 2/18/2014, 8:30-9:05 AM) ---
# B6/C6/D6 inherited row 5's special (left-aligned date/time) formatting from
# the Insert() above; clear that back off so they use the same plain
# General/h:mm formatting as the other data rows (e.g. row 4) use.
$ws.Cells.Item(6, 2).ClearFormats()
$ws.Cells.Item(6, 3).ClearFormats()
$ws.Cells.Item(6, 3).NumberFormat = "h:mm"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "h:mm"

$ws.Cells.Item(6, 1).Value = 41688
$ws.Cells.Item(6, 2).Value = "Doyle"
$ws.Cells.Item(6, 3).Value = 0.35416666666666669
$ws.Cells.Item(6, 4).Value = 0.37847222222222227
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 9).Formula = "=35"
$ws.Cells.Item(6, 11).Value = "Put together start of data definitions (world, mouse and key events)"

# --- Row 7: mirrors row 5's "committing to git" summary formulas, now referencing row 6 ---
$ws.Cells.Item(7, 1).Value = " =========================    committing to git:"
$ws.Cells.Item(7, 2).Formula = "=A6"
$ws.Cells.Item(7, 3).Formula = "=D6"
$ws.Cells.Item(7, 4).Value = " =========================    "

# Fix up the selected cell to match the author's final position
$ws.Range("K6").Select()
